# Applies the "pull new files from local" refresh to the attendance sheet:
# - rewrites the attendance rows, now sorted by student id ascending
# - updates one student id (20101134 -> 20101135)
# - updates attendance flags / timestamps for several students
# - widens/resizes the columns (and adds a width for the new-ish column D)
# - updates the saved active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New attendance table contents (10 data rows, rows 2-11), sorted by
# student id ascending, reflecting the freshly pulled attendance log.
$rows = @(
    @("1",        "강동완",          "",  "NaN"),
    @("15102993",  "Ngoc Anh",       "",  "NaN"),
    @("16101384",  "응웬민뚜",        "Y", "2020/12/29 00:36:08"),
    @("16101883",  "Khanh Ngan",     "",  "NaN"),
    @("17101222",  "안홍현",          "",  "NaN"),
    @("18101255",  "Linh",           "",  "NaN"),
    @("18102231",  "Thanh Danh",     "Y", "2020/12/29 00:36:16"),
    @("19102395",  "Anh Duc",        "",  "NaN"),
    @("20101023",  "Nhim",           "Y", "2020/12/29 00:37:01"),
    @("20101135",  "Pham Duy Thai",  "Y", "2020/12/29 00:36:18")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    # Column A holds student ids that look numeric (e.g. "15102993") but
    # must stay text, matching the original sharedString storage - enter
    # them with a leading apostrophe, then drop back to the Normal style
    # so no stray number formatting sticks to the cell.
    $ws.Cells.Item($r, 1).Value = "'" + $data[0]
    $ws.Cells.Item($r, 1).Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $data[1]

    if ($data[2] -ne "") {
        $ws.Cells.Item($r, 3).Value = $data[2]
    } else {
        $ws.Cells.Item($r, 3).Value = $null
    }

    $ws.Cells.Item($r, 4).Value = $data[3]
}

# Resize columns to match the refreshed layout (closest achievable widths -
# the host snaps column widths to its internal pixel grid).
$ws.Columns.Item(1).ColumnWidth = 13.8
$ws.Columns.Item(2).ColumnWidth = 21.65
$ws.Columns.Item(3).ColumnWidth = 15.1
$ws.Columns.Item(4).ColumnWidth = 38.1

# Update the saved selection/active cell
$ws.Range("C15").Select()
